# Weekly update: insert a new daily price record as the new first data
# row (row 172) for Choclo / Comercializadora del Agro de Limarí,
# pushing the existing rows 172:200 down to 173:201.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 172; this shifts rows
# 172:200 down to 173:201 and extends the sheet to 201 rows.
$ws.Rows("172").Insert()

# Populate the new row 172 with the latest observation.
$ws.Range("A172").Value = 2
$ws.Range("B172").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C172").Value = "Coquimbo"
$ws.Range("D172").Value = 45211
$ws.Range("E172").Value = 4
$ws.Range("F172").Value = 100112024
$ws.Range("G172").Value = "Choclo"
$ws.Range("H172").Value = "Dulce o Americano"
$ws.Range("I172").Value = "Primera"
$ws.Range("J172").Value = 700
$ws.Range("K172").Value = 38000
$ws.Range("L172").Value = 40000
$ws.Range("M172").Value = 39000
$ws.Range("N172").Value = "`$/malla 70 unidades"
$ws.Range("O172").Value = "Provincia de Limarí"
$ws.Range("P172").Value = 557
$ws.Range("Q172").Value = 70
$ws.Range("R172").Value = "Hortaliza"
